$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("input")

# New row 9: "additional_accounts" record (covid and other additions)
$ws.Range("A9").Value = "additional_accounts"

# Column B holds the same literal "1.0" text as every other row; copy it from
# B8 so the new cell picks up the identical shared-string/text typing instead
# of being reinterpreted as a number.
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4104)  # xlPasteAll

$ws.Range("C9:R9").Value = 0

# S9 (2021) uses the same highlighted thousands-style as the other rows' S column.
$ws.Range("S8").Copy()
$ws.Range("S9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("S9").Value = 11920

$excel.CutCopyMode = 0

# Move the active selection, matching the author's final cursor position.
$ws.Range("M5").Select()
